# Update the cryptos worksheet with the latest prices/volumes scraped by GitHub Actions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.766.57'
$ws.Range('E2').Value = '  +4.50%  '
$ws.Range('D3').Value = '2.500.90'
$ws.Range('E3').Value = '  +2.34%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''489.00'
$ws.Range('E5').Value = '  +4.39%  '
$ws.Range('D6').Value = '''146.86'
$ws.Range('E6').Value = '  +10.00%  '
$ws.Range('D7').Value = '''0.998'
$ws.Range('E7').Value = '  +0.34%  '
$ws.Range('D8').Value = '''0.514'
$ws.Range('E8').Value = '  +4.90%  '
$ws.Range('D9').Value = '2.521.53'
$ws.Range('E9').Value = '  +2.30%  '
$ws.Range('D10').Value = '''5.79'
$ws.Range('E10').Value = '  +8.24%  '
$ws.Range('D11').Value = '''0.0978'
$ws.Range('E11').Value = '  +1.55%  '
$ws.Range('D12').Value = '''0.332'
$ws.Range('E12').Value = '  +4.12%  '
$ws.Range('E13').Value = '  +1.17%  '
$ws.Range('D14').Value = '2.944.53'
$ws.Range('E14').Value = '  +3.07%  '
$ws.Range('D15').Value = '56.635.56'
$ws.Range('E15').Value = '  +4.56%  '
$ws.Range('D16').Value = '''21.27'
$ws.Range('E16').Value = '  +6.78%  '
$ws.Range('D17').Value = '''0.0000137'
$ws.Range('E17').Value = '  +1.93%  '
$ws.Range('D18').Value = '2.522.75'
$ws.Range('E18').Value = '  +3.03%  '
$ws.Range('D19').Value = '''4.52'
$ws.Range('E19').Value = '  +6.99%  '
$ws.Range('D20').Value = '''10.22'
$ws.Range('E20').Value = '  +7.90%  '
$ws.Range('D21').Value = '''321.46'
$ws.Range('E21').Value = '  +2.07%  '
$ws.Range('D22').Value = '''0.999'
$ws.Range('E22').Value = '  +0.70%  '
$ws.Range('D23').Value = '''5.84'
$ws.Range('E23').Value = '  +7.93%  '
$ws.Range('D24').Value = '''58.84'
$ws.Range('E24').Value = '  +3.58%  '
$ws.Range('D25').Value = '''0.413'
$ws.Range('E25').Value = '  +6.64%  '
$ws.Range('D26').Value = '''0.167'
$ws.Range('E26').Value = '  +8.55%  '
$ws.Range('D27').Value = '''0.998'
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('D28').Value = '2.621.84'
$ws.Range('E28').Value = '  +3.19%  '
$ws.Range('D29').Value = '''7.65'
$ws.Range('E29').Value = '  +5.09%  '
$ws.Range('D30').Value = '0.0₃0798'
$ws.Range('E30').Value = '  +9.04%  '
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('D32').Value = '''148.97'
$ws.Range('E32').Value = '  -0.83%  '
$ws.Range('D33').Value = '''18.28'
$ws.Range('E33').Value = '  +2.57%  '
$ws.Range('E34').Value = '  +5.01%  '
$ws.Range('D35').Value = '''5.21'
$ws.Range('E35').Value = '  +3.16%  '
$ws.Range('E36').Value = '  +7.65%  '
$ws.Range('D37').Value = '''3.75'
$ws.Range('E37').Value = '  +4.44%  '
$ws.Range('D38').Value = '''0.871'
$ws.Range('E38').Value = '  +7.63%  '
$ws.Range('D39').Value = '''34.30'
$ws.Range('E39').Value = '  +1.78%  '
$ws.Range('D40').Value = '''3.55'
$ws.Range('E40').Value = '  +7.52%  '
$ws.Range('D41').Value = '''0.619'
$ws.Range('E41').Value = '  +2.29%  '
$ws.Range('D42').Value = '''0.0559'
$ws.Range('E42').Value = '  +5.08%  '
$ws.Range('D43').Value = '''0.996'
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('E44').Value = '  +7.28%  '
$ws.Range('D45').Value = '''4.81'
$ws.Range('E45').Value = '  +10.29%  '
$ws.Range('D46').Value = '''261.54'
$ws.Range('E46').Value = '  +17.59%  '
$ws.Range('D47').Value = '''0.0229'
$ws.Range('E47').Value = '  +3.39%  '
$ws.Range('B48').Value = 'WhiteBITCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D48').Value = '''10.20'
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '''0.0910'
$ws.Range('E49').Value = '  +4.04%  '
$ws.Range('D50').Value = '1.916.31'
$ws.Range('E50').Value = '  -2.20%  '
$ws.Range('D51').Value = '''17.72'
$ws.Range('E51').Value = '  +6.18%  '
